# Weights.xlsx — add the two missing weight values (Fuselage 1 / Rudder
# Bottom) on Sheet1 and move the saved cursor/scroll position, matching
# what the workbook's author did before re-uploading the file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weight entries
$ws.Range("B3").Value = 11
$ws.Range("B18").Value = 2

# Move the saved selection to B18 and scroll the sheet down so row 10
# is at the top of the visible area.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B18").Select()
